$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regression")

# New columns F and G widths (best-fit to match "hotel"/"Hotel Creek" and
# "Check In Date" contents; the runtime quantizes ColumnWidth to 1/6-character
# steps, so these are the closest achievable values to the recorded
# 10.6640625 / 12.33203125 XML widths)
$ws.Columns.Item(6).ColumnWidth = 9.75
$ws.Columns.Item(7).ColumnWidth = 11.42

# Fill the new block column-by-column (A7,A8,B7,B8,...) to mirror the
# order the values were originally entered in Excel.
$ws.Range("A7").Value = "TC-113"
$ws.Range("A8").Value = "TC-113"

$ws.Range("B7").Value = "username"
$ws.Range("B8").Value = "reyaz0806"

$ws.Range("C7").Value = "password"
$ws.Range("C8").Value = "reyaz123"

$ws.Range("D7").Value = "expected Title"
$ws.Range("D8").Value = "Adactin.com - Search Hotel"

$ws.Range("E7").Value = "location"
$ws.Range("E8").Value = "Sydney"

$ws.Range("F7").Value = "hotel"
$ws.Range("F8").Value = "Hotel Creek"

$ws.Range("G7").Value = "Check In Date"
$ws.Range("G8").Value = "'27/09/2024"
$ws.Range("G8").NumberFormat = "mm-dd-yy"

# Update the selection to mimic the saved cursor position
$ws.Range("D10").Select()
